$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 2485.9285
$ws.Range("I115").Value = 184.6
$ws.Range("J115").Value = 3764.4443
$ws.Range("K115").Value = 553.8
$ws.Range("L115").Value = 11293.3329
$ws.Range("M115").Value = 1013.2
$ws.Range("N115").Value = -14427.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26992.354
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 26992.354
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 26992.354
$ws.Range("N32").Value = -27566.354
$ws.Range("M32").ClearContents()

$ws.Range("H61").Value = 4214.4326
$ws.Range("I61").Value = 4783.5356
$ws.Range("J61").Value = 2443.889
$ws.Range("K61").Value = 4783.5356
$ws.Range("L61").Value = 2443.889
$ws.Range("M61").Value = -4571.5356
$ws.Range("N61").Value = -2867.889

$ws.Range("H80").Value = 30076.666
$ws.Range("J80").Value = 30076.666
$ws.Range("L80").Value = 30076.666
$ws.Range("N80").Value = -32072.666

$ws.Range("H83").Value = 30076.666
$ws.Range("J83").Value = 30076.666
$ws.Range("L83").Value = 90229.99800000001
$ws.Range("N83").Value = -100213.998

$ws.Range("H110").Value = 1100
$ws.Range("I110").Value = 1100
$ws.Range("K110").Value = 1100
$ws.Range("M110").Value = 945

$ws.Range("H134").Value = 21036.363
$ws.Range("J134").Value = 21036.363
$ws.Range("L134").Value = 21036.363
$ws.Range("N134").Value = -31176.363

$ws.Range("H136").Value = 4214.4326
$ws.Range("I136").Value = 4783.5356
$ws.Range("J136").Value = 2443.889
$ws.Range("K136").Value = 14350.6068
$ws.Range("L136").Value = 7331.667
$ws.Range("M136").Value = -11800.6068
$ws.Range("N136").Value = -12431.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 143.25
$ws.Range("I35").Value = 143.25
$ws.Range("K35").Value = 143.25
$ws.Range("M35").Value = 150.75

$ws.Range("H52").Value = 39003.89
$ws.Range("I52").Value = 10180
$ws.Range("J52").Value = 47239.285
$ws.Range("K52").Value = 10180
$ws.Range("L52").Value = 47239.285
$ws.Range("M52").Value = -9886
$ws.Range("N52").Value = -47827.285

$ws.Range("H134").Value = 52502250
$ws.Range("I134").Value = 125002330
$ws.Range("J134").Value = 4168860
$ws.Range("K134").Value = 375006990
$ws.Range("L134").Value = 12506580
$ws.Range("M134").Value = -375004455
$ws.Range("N134").Value = -12511650

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 436.6111
$ws.Range("I4").Value = 110
$ws.Range("J4").Value = 844.875
$ws.Range("K4").Value = 330
$ws.Range("L4").Value = 2534.625
$ws.Range("M4").Value = -218
$ws.Range("N4").Value = -2758.625

$ws.Range("H75").Value = 1670.909
$ws.Range("J75").Value = 1711.1111
$ws.Range("L75").Value = 5133.3333
$ws.Range("N75").Value = -7129.3333

$ws.Range("H76").Value = 200
$ws.Range("I76").Value = 200
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 600
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -217
$ws.Range("N76").ClearContents()

$ws.Range("H78").Value = 1670.909
$ws.Range("J78").Value = 1711.1111
$ws.Range("L78").Value = 15399.9999
$ws.Range("N78").Value = -25383.9999

$ws.Range("H79").Value = 200
$ws.Range("I79").Value = 200
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 600
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 726
$ws.Range("N79").ClearContents()

$ws.Range("H94").Value = 3279.5715
$ws.Range("I94").Value = 1066.6666
$ws.Range("J94").Value = 3487.0312
$ws.Range("K94").Value = 3199.9998
$ws.Range("L94").Value = 10461.0936
$ws.Range("M94").Value = -2523.9998
$ws.Range("N94").Value = -11813.0936

$ws.Range("H118").Value = 76924264
$ws.Range("I118").Value = 100000790
$ws.Range("J118").Value = 2500
$ws.Range("K118").Value = 300002370
$ws.Range("L118").Value = 7500
$ws.Range("M118").Value = -300001127
$ws.Range("N118").Value = -9986

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 295.55554
$ws.Range("I107").Value = 93.85714
$ws.Range("J107").Value = 1001.5
$ws.Range("K107").Value = 93.85714
$ws.Range("L107").Value = 1001.5
$ws.Range("M107").Value = 1826.14286
$ws.Range("N107").Value = -4841.5

$ws.Range("H141").Value = 23047.666
$ws.Range("J141").Value = 23047.666
$ws.Range("L141").Value = 23047.666
$ws.Range("N141").Value = -33407.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1123.8823
$ws.Range("I7").Value = 1157.875
$ws.Range("J7").Value = 580
$ws.Range("K7").Value = 1157.875
$ws.Range("L7").Value = 580
$ws.Range("M7").Value = -1045.875
$ws.Range("N7").Value = -804

$ws.Range("H9").Value = 3985
$ws.Range("I9").Value = 272.66666
$ws.Range("J9").Value = 7697.3335
$ws.Range("K9").Value = 272.66666
$ws.Range("L9").Value = 7697.3335
$ws.Range("M9").Value = -48.66665999999998
$ws.Range("N9").Value = -8145.3335

$ws.Range("H35").Value = 400
$ws.Range("I35").Value = 550
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 550
$ws.Range("L35").Value = 100
$ws.Range("M35").Value = -214
$ws.Range("N35").Value = -772

$ws.Range("H100").Value = 2831.6667
$ws.Range("I100").Value = 980
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 980
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -439
$ws.Range("N100").Value = -4082

$ws.Range("H126").Value = 1123.8823
$ws.Range("I126").Value = 1157.875
$ws.Range("J126").Value = 580
$ws.Range("K126").Value = 3473.625
$ws.Range("L126").Value = 1740
$ws.Range("M126").Value = -1003.625
$ws.Range("N126").Value = -6680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64904.668
$ws.Range("J46").Value = 64904.668
$ws.Range("L46").Value = 64904.668
$ws.Range("N46").Value = -65366.668

$ws.Range("H81").Value = 1984.3334
$ws.Range("I81").Value = 678.5
$ws.Range("J81").Value = 2506.6667
$ws.Range("K81").Value = 1357
$ws.Range("L81").Value = 5013.3334
$ws.Range("M81").Value = -296
$ws.Range("N81").Value = -7135.3334

$ws.Range("H84").Value = 1984.3334
$ws.Range("I84").Value = 678.5
$ws.Range("J84").Value = 2506.6667
$ws.Range("K84").Value = 6785
$ws.Range("L84").Value = 25066.667
$ws.Range("M84").Value = -1481
$ws.Range("N84").Value = -35674.667

$ws.Range("H113").Value = 583.5
$ws.Range("I113").Value = 471.44446
$ws.Range("J113").Value = 1188.6
$ws.Range("K113").Value = 1414.33338
$ws.Range("L113").Value = 3565.8
$ws.Range("M113").Value = 755.66662
$ws.Range("N113").Value = -7905.799999999999

$ws.Range("H125").Value = 48358
$ws.Range("J125").Value = 48358
$ws.Range("L125").Value = 48358
$ws.Range("N125").Value = -58198

$ws.Range("H134").Value = 64904.668
$ws.Range("J134").Value = 64904.668
$ws.Range("L134").Value = 194714.004
$ws.Range("N134").Value = -199784.004
